# Atualização de bases das ligas, do dia: 14-05-2024 às 01:09
#
# The source data for several fixtures had their rows swapped with an
# adjacent fixture row (same kick-off date/time). Column A (the row's
# positional rank index) stays put; every other column (B..AB — id,
# HomeTeam, AwayTeam, score, odds, etc.) moves with its fixture, so we
# swap the full B:AB payload between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B${rowA}:AB${rowA}")
    $rangeB = $ws.Range("B${rowB}:AB${rowB}")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 47 48
Swap-Rows 104 105
Swap-Rows 107 108
Swap-Rows 128 129
Swap-Rows 142 145
Swap-Rows 153 154
